# This script re-applies a re-scrape of the "electronics" product listing:
#   1) A handful of rows were re-ordered (the underlying scrape re-sorted /
#      re-paginated slightly differently between runs), so their entire
#      row content (columns A-N) needs to move to the new row positions.
#   2) Every data row's timestamp (column O) is bumped to the new scrape
#      time "2022-08-23 20:57:41".
#
# Columns: A id, B title, C href, D quantity, E ratingAmount, F ratingValue,
#          G brand, H price, I priceContext, J priceContextHiddenText,
#          K priceContextPrice, L priceContextAmount, M udoCat,
#          N productAriaLabel, O timestamp
# Columns E/F hold real numbers; every other column (even the numeric
# looking ones like H/K) is stored as text, so we force a text
# NumberFormat before writing those so Excel doesn't silently convert
# "9.95" into the number 9.95.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowValues($sheet, $row) {
    $vals = @()
    for ($c = 1; $c -le 15; $c++) {
        $vals += ,$sheet.Cells.Item($row, $c).Value2
    }
    return $vals
}

function Set-RowValues($sheet, $row, $vals) {
    for ($c = 1; $c -le 15; $c++) {
        $cell = $sheet.Cells.Item($row, $c)
        if ($c -eq 5 -or $c -eq 6) {
            # E (ratingAmount) / F (ratingValue) are real numbers
            if ($vals[$c - 1] -eq $null -or $vals[$c - 1] -eq "") {
                $cell.Value = ""
            } else {
                $cell.Value = $vals[$c - 1]
            }
        } else {
            $cell.NumberFormat = "@"
            $cell.Value = [string]$vals[$c - 1]
        }
    }
}

# New-row -> old-row source mapping for every row whose content moved.
$rowMap = @{
    18 = 19; 19 = 18;
    29 = 31; 31 = 29;
    38 = 39; 39 = 38;
    42 = 46; 43 = 42; 46 = 50; 47 = 43; 48 = 47; 49 = 48; 50 = 49;
    53 = 55; 54 = 53; 55 = 54;
    60 = 63; 61 = 60; 63 = 61;
    64 = 66; 65 = 64; 66 = 65;
}

# Snapshot every involved row BEFORE writing anything, since several rows
# take part in multi-row cycles (e.g. 42 -> 46 -> 50 -> 49 -> 48 -> 47 -> 43
# -> 42) and would otherwise clobber each other mid-update.
$snapshot = @{}
foreach ($row in $rowMap.Keys) {
    $snapshot[$row] = Get-RowValues $ws $row
}

foreach ($row in $rowMap.Keys) {
    $sourceRow = $rowMap[$row]
    Set-RowValues $ws $row $snapshot[$sourceRow]
}

# Bump the scrape timestamp on every data row (2-86).
$tsRange = $ws.Range("O2:O86")
$tsRange.NumberFormat = "@"
$tsRange.Value = "2022-08-23 20:57:41"
